$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump the generation Date stamp
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-30T16:36:55+00:00"

# ---------------------------------------------------------------------------
# 2. Elements sheet: widen column K (Type(s)) and append a new element row
#    (VieQuotidienne.PersonnePriseCharge) that mirrors the layout of row 5.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Widen column 11 ("Type(s)") to fit the new long URL value.
$ws.Columns.Item(11).ColumnWidth = 61.8

# Clone formatting of the last data row (row 5) into the new row 6 first, so
# every cell in row 6 starts out with the same style (s="2") as the rest of
# the table before any values are written into it.
$ws.Range("A5:AJ5").Copy()
$ws.Range("A6:AJ6").PasteSpecial(-4122)

# --- plain text values -----------------------------------------------------
$ws.Range("A6").Value = "VieQuotidienne.PersonnePriseCharge"
$ws.Range("B6").Value = "VieQuotidienne.PersonnePriseCharge"
$ws.Range("K6").Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/PersonnePriseCharge`n"
$ws.Range("L6").Value = "Lien vers la classe PersonnePriseCharge"
$ws.Range("M6").Value = "Lien vers la classe PersonnePriseCharge"
$ws.Range("AF6").Value = "VieQuotidienne.PersonnePriseCharge"

# --- numeric-looking text values ("1") and blank-but-typed text ("") ------
# These columns hold their content as *text*, not numbers, in the source
# workbook, so a leading apostrophe is used to force text entry. The
# PasteSpecial(formats) pass afterwards strips the resulting quote-prefix
# styling back to the shared row style (s="2").
$textCells = @("D6","F6","G6","H6","I6","J6","P6","R6","S6","T6","U6","V6","W6","X6","Y6","Z6","AA6","AB6","AC6","AD6","AE6","AG6","AH6","AI6","AJ6")
$textValues = @{
    "D6" = "'"
    "F6" = "'1"
    "G6" = "'1"
    "H6" = "'"
    "I6" = "'"
    "J6" = "'"
    "P6" = "'"
    "R6" = "'"
    "S6" = "'"
    "T6" = "'"
    "U6" = "'"
    "V6" = "'"
    "W6" = "'"
    "X6" = "'"
    "Y6" = "'"
    "Z6" = "'"
    "AA6" = "'"
    "AB6" = "'"
    "AC6" = "'"
    "AD6" = "'"
    "AE6" = "'"
    "AG6" = "'1"
    "AH6" = "'1"
    "AI6" = "'"
    "AJ6" = "'"
}

foreach ($ref in $textCells) {
    $ws.Range($ref).Value = $textValues[$ref]
}

# Re-paste the original (row 5) formatting onto those same cells so the
# quote-prefix style introduced above is replaced by the shared body style.
foreach ($ref in $textCells) {
    $col = $ref -replace '6$', ''
    $ws.Range("$col`5").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

Write-Output "done"
